$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1815.5834
$ws.Range("I40").Value = 1590
$ws.Range("J40").Value = 1928.375
$ws.Range("K40").Value = 1590
$ws.Range("L40").Value = 1928.375
$ws.Range("M40").Value = -1415
$ws.Range("N40").Value = -2278.375

$ws.Range("H80").Value = 8122531.5
$ws.Range("I80").Value = 604.4666999999999
$ws.Range("K80").Value = 1813.4001
$ws.Range("M80").Value = -815.4000999999998

$ws.Range("H83").Value = 8122531.5
$ws.Range("I83").Value = 604.4666999999999
$ws.Range("K83").Value = 5440.2003
$ws.Range("M83").Value = -448.2002999999995

$ws.Range("H116").Value = 4549.6
$ws.Range("I116").Value = 1500
$ws.Range("J116").Value = 6582.6665
$ws.Range("K116").Value = 1500
$ws.Range("L116").Value = 6582.6665
$ws.Range("M116").Value = 1942
$ws.Range("N116").Value = -13466.6665

$ws.Range("H138").Value = 1903.4675
$ws.Range("J138").Value = 2243.3276
$ws.Range("L138").Value = 6729.9828
$ws.Range("N138").Value = -17009.9828

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22865.672
$ws.Range("I32").Value = 26683.674
$ws.Range("K32").Value = 26683.674
$ws.Range("M32").Value = -26396.674

$ws.Range("H74").Value = 40001100
$ws.Range("I74").Value = 62500580
$ws.Range("J74").Value = 2022.7778
$ws.Range("K74").Value = 62500580
$ws.Range("L74").Value = 2022.7778
$ws.Range("M74").Value = -62499706
$ws.Range("N74").Value = -3770.7778

$ws.Range("H77").Value = 40001100
$ws.Range("I77").Value = 62500580
$ws.Range("J77").Value = 2022.7778
$ws.Range("K77").Value = 312502900
$ws.Range("L77").Value = 10113.889
$ws.Range("M77").Value = -312498532
$ws.Range("N77").Value = -18849.889

$ws.Range("H102").Value = 1585
$ws.Range("I102").Value = 1255
$ws.Range("K102").Value = 1255
$ws.Range("M102").Value = 367

$ws.Range("H122").Value = 2797
$ws.Range("I122").Value = 1740.9
$ws.Range("K122").Value = 5222.700000000001
$ws.Range("M122").Value = -2772.700000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9577.232
$ws.Range("I31").Value = 17158.6
$ws.Range("K31").Value = 17158.6
$ws.Range("M31").Value = -16863.6

$ws.Range("H34").Value = 9577.232
$ws.Range("I34").Value = 17158.6
$ws.Range("K34").Value = 17158.6
$ws.Range("M34").Value = -16956.6

$ws.Range("H58").Value = 16080.637
$ws.Range("I58").Value = 835.2917
$ws.Range("J58").Value = 56734.89
$ws.Range("K58").Value = 835.2917
$ws.Range("L58").Value = 56734.89
$ws.Range("M58").Value = -632.2917
$ws.Range("N58").Value = -57140.89

$ws.Range("H136").Value = 16080.637
$ws.Range("I136").Value = 835.2917
$ws.Range("J136").Value = 56734.89
$ws.Range("K136").Value = 2505.8751
$ws.Range("L136").Value = 170204.67
$ws.Range("M136").Value = 44.1248999999998
$ws.Range("N136").Value = -175304.67

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1205.5714
$ws.Range("J23").Value = 1336.1818
$ws.Range("L23").Value = 4008.5454
$ws.Range("N23").Value = -4478.5454

$ws.Range("H58").Value = 3720.6667
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H132").Value = 2044.4445
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2044.4445
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 18400.0005
$ws.Range("N132").Value = -23460.0005
$ws.Range("M132").ClearContents()

$ws.Range("H137").Value = 23812234
$ws.Range("J137").Value = 37040544
$ws.Range("L137").Value = 111121632
$ws.Range("N137").Value = -111131832

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 30000
$ws.Range("J6").Value = 30000
$ws.Range("L6").Value = 30000
$ws.Range("N6").Value = -30226

$ws.Range("H16").Value = 30000
$ws.Range("J16").Value = 30000
$ws.Range("L16").Value = 30000
$ws.Range("N16").Value = -30500

$ws.Range("H70").Value = 2023319.6
$ws.Range("I70").Value = 9195.210999999999
$ws.Range("J70").Value = 5212350
$ws.Range("K70").Value = 9195.210999999999
$ws.Range("L70").Value = 5212350
$ws.Range("M70").Value = -8925.210999999999
$ws.Range("N70").Value = -5212890

$ws.Range("H73").Value = 2023319.6
$ws.Range("I73").Value = 9195.210999999999
$ws.Range("J73").Value = 5212350
$ws.Range("K73").Value = 9195.210999999999
$ws.Range("L73").Value = 5212350
$ws.Range("M73").Value = -8259.210999999999
$ws.Range("N73").Value = -5214222

$ws.Range("H102").Value = 33334696
$ws.Range("I102").Value = 38462844
$ws.Range("J102").Value = 1749.5
$ws.Range("K102").Value = 38462844
$ws.Range("L102").Value = 1749.5
$ws.Range("M102").Value = -38461222
$ws.Range("N102").Value = -4993.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4431.9165
$ws.Range("I7").Value = 2736.8
$ws.Range("J7").Value = 5642.7144
$ws.Range("K7").Value = 2736.8
$ws.Range("L7").Value = 5642.7144
$ws.Range("M7").Value = -2624.8
$ws.Range("N7").Value = -5866.7144

$ws.Range("H16").Value = 755.5714
$ws.Range("I16").Value = 755.5714
$ws.Range("K16").Value = 755.5714
$ws.Range("M16").Value = -585.5714

$ws.Range("H61").Value = 4835.88
$ws.Range("J61").Value = 7353.3
$ws.Range("L61").Value = 7353.3
$ws.Range("N61").Value = -7757.3

$ws.Range("H113").Value = 4835.88
$ws.Range("J113").Value = 7353.3
$ws.Range("L113").Value = 7353.3
$ws.Range("N113").Value = -11693.3

$ws.Range("H126").Value = 4431.9165
$ws.Range("I126").Value = 2736.8
$ws.Range("J126").Value = 5642.7144
$ws.Range("K126").Value = 8210.400000000001
$ws.Range("L126").Value = 16928.1432
$ws.Range("M126").Value = -5740.400000000001
$ws.Range("N126").Value = -21868.1432

$ws.Range("H132").Value = 2335.0588
$ws.Range("I132").Value = 1672.5454
$ws.Range("K132").Value = 5017.6362
$ws.Range("M132").Value = -2487.6362

$ws.Range("H136").Value = 25837.75
$ws.Range("I136").Value = 34020.668
$ws.Range("J136").Value = 1289
$ws.Range("K136").Value = 102062.004
$ws.Range("L136").Value = 3867
$ws.Range("M136").Value = -99512.00399999999
$ws.Range("N136").Value = -8967

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 36316.668
$ws.Range("I18").Value = 2950
$ws.Range("J18").Value = 53000
$ws.Range("K18").Value = 2950
$ws.Range("L18").Value = 53000
$ws.Range("M18").Value = -2777
$ws.Range("N18").Value = -53346

$ws.Range("H122").Value = 2232.8333
$ws.Range("I122").Value = 2079.6
$ws.Range("K122").Value = 6238.799999999999
$ws.Range("M122").Value = -3788.799999999999

$ws.Range("H132").Value = 1149.1111
$ws.Range("I132").Value = 656.8570999999999
$ws.Range("J132").Value = 2872
$ws.Range("K132").Value = 1970.5713
$ws.Range("L132").Value = 8616
$ws.Range("M132").Value = 559.4287000000002
$ws.Range("N132").Value = -13676
